$wb = $excel.ActiveWorkbook

# Row 17 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1227.5
$ws.Range("J17").Value = 1227.5
$ws.Range("L17").Value = 3682.5
$ws.Range("N17").Value = -4018.5

# Row 98 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4278.409
$ws.Range("I98").Value = 4425
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 4425
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = -2927
$ws.Range("N98").Value = -4196

# Row 100 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1615.1852
$ws.Range("I100").Value = 1423.1333
$ws.Range("J100").Value = 1855.25
$ws.Range("K100").Value = 1423.1333
$ws.Range("L100").Value = 1855.25
$ws.Range("M100").Value = -882.1333
$ws.Range("N100").Value = -2937.25

# Row 116 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2979.8
$ws.Range("I116").Value = 2950
$ws.Range("J116").Value = 2999.6667
$ws.Range("K116").Value = 2950
$ws.Range("L116").Value = 2999.6667
$ws.Range("M116").Value = 492
$ws.Range("N116").Value = -9883.6667

# Row 118 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1098
$ws.Range("I118").Value = 747.5
$ws.Range("J118").Value = 2500
$ws.Range("K118").Value = 2242.5
$ws.Range("L118").Value = 7500
$ws.Range("M118").Value = -585.5
$ws.Range("N118").Value = -10814

# Row 122 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4278.409
$ws.Range("I122").Value = 4425
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 13275
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -10825
$ws.Range("N122").Value = -8500

# Row 132 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7460
$ws.Range("I132").Value = 7707.143
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 23121.429
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -20591.429
$ws.Range("N132").Value = -17060

# Row 61 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2483.7307
$ws.Range("I61").Value = 2717.3125
$ws.Range("J61").Value = 2110
$ws.Range("K61").Value = 2717.3125
$ws.Range("L61").Value = 2110
$ws.Range("M61").Value = -2505.3125
$ws.Range("N61").Value = -2534

# Row 69 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 79999
$ws.Range("J69").Value = 79999
$ws.Range("L69").Value = 79999
$ws.Range("N69").Value = -81497

# Row 72 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H72").Value = 79999
$ws.Range("J72").Value = 79999
$ws.Range("L72").Value = 239997
$ws.Range("N72").Value = -247485

# Row 110 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1504.8235
$ws.Range("I110").Value = 1262.75
$ws.Range("J110").Value = 1720
$ws.Range("K110").Value = 1262.75
$ws.Range("L110").Value = 1720
$ws.Range("M110").Value = 782.25
$ws.Range("N110").Value = -5810

# Row 132 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 30168.611
$ws.Range("I132").Value = 2334.2068
$ws.Range("J132").Value = 145482.58
$ws.Range("K132").Value = 7002.6204
$ws.Range("L132").Value = 436447.74
$ws.Range("M132").Value = -4472.6204
$ws.Range("N132").Value = -441507.74

# Row 136 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2483.7307
$ws.Range("I136").Value = 2717.3125
$ws.Range("J136").Value = 2110
$ws.Range("K136").Value = 8151.9375
$ws.Range("L136").Value = 6330
$ws.Range("M136").Value = -5601.9375
$ws.Range("N136").Value = -11430

# Row 94 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 457.55554
$ws.Range("I94").Value = 273.08334
$ws.Range("K94").Value = 273.08334
$ws.Range("M94").Value = 177.91666

# Row 99 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 907.8946999999999
$ws.Range("I99").Value = 709.8461
$ws.Range("K99").Value = 709.8461
$ws.Range("M99").Value = 788.1539

# Row 113 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 613.5599999999999
$ws.Range("I113").Value = 668.75
$ws.Range("J113").Value = 587.58826
$ws.Range("K113").Value = 2006.25
$ws.Range("L113").Value = 1762.76478
$ws.Range("M113").Value = 163.75
$ws.Range("N113").Value = -6102.76478

# Row 121 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I121").Value = 532.1111
$ws.Range("K121").Value = 1596.3333
$ws.Range("M121").Value = -286.3332999999998

# Row 131 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2179.7158
$ws.Range("I131").Value = 10426.667
$ws.Range("J131").Value = 1576.2805
$ws.Range("K131").Value = 31280.001
$ws.Range("L131").Value = 4728.8415
$ws.Range("M131").Value = -26240.001
$ws.Range("N131").Value = -14808.8415

# Row 16 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2878.4666
$ws.Range("I16").Value = 2726.9285
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 2726.9285
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -2556.9285
$ws.Range("N16").Value = -5340

# Row 61 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1429
$ws.Range("I61").Value = 1349.3572
$ws.Range("J61").Value = 1652
$ws.Range("K61").Value = 1349.3572
$ws.Range("L61").Value = 1652
$ws.Range("M61").Value = -1147.3572
$ws.Range("N61").Value = -2056

# Row 113 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1429
$ws.Range("I113").Value = 1349.3572
$ws.Range("J113").Value = 1652
$ws.Range("K113").Value = 1349.3572
$ws.Range("L113").Value = 1652
$ws.Range("M113").Value = 820.6428000000001
$ws.Range("N113").Value = -5992

# Row 122 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2837.5
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 3280
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 9840
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -14740

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5383.6587
$ws.Range("I132").Value = 5467.5557
$ws.Range("J132").Value = 4779.6
$ws.Range("K132").Value = 16402.6671
$ws.Range("L132").Value = 14338.8
$ws.Range("M132").Value = -13872.6671
$ws.Range("N132").Value = -19398.8

# Row 136 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1970.8695
$ws.Range("I136").Value = 1195.8823
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 3587.6469
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -1037.6469
$ws.Range("N136").Value = -17599.9995

# Row 107 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 625.3684
$ws.Range("I107").Value = 529.2222
$ws.Range("J107").Value = 711.9
$ws.Range("K107").Value = 1587.6666
$ws.Range("L107").Value = 2135.7
$ws.Range("M107").Value = 332.3334
$ws.Range("N107").Value = -5975.7

# Row 113 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 523.5833
$ws.Range("I113").Value = 542.55554
$ws.Range("J113").Value = 466.66666
$ws.Range("K113").Value = 1627.66662
$ws.Range("L113").Value = 1399.99998
$ws.Range("M113").Value = 542.33338
$ws.Range("N113").Value = -5739.999980000001

# Row 122 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2179.7917
$ws.Range("I122").Value = 1561.9445
$ws.Range("J122").Value = 4033.3333
$ws.Range("K122").Value = 4685.833500000001
$ws.Range("L122").Value = 12099.9999
$ws.Range("M122").Value = -2235.833500000001
$ws.Range("N122").Value = -16999.9999

# Row 126 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2207.0312
$ws.Range("I126").Value = 1944.2307
$ws.Range("K126").Value = 5832.6921
$ws.Range("M126").Value = -3362.6921
